$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: copy format from E1 (same header style), then set its value/label
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# time_taken values for rows 2..64 (plain, unstyled cells like the rest of the data rows)
$timeTaken = @(
    "2021-10-05 13:40:48.315668",
    "2021-10-05 13:40:48.315680",
    "2021-10-05 13:40:48.315683",
    "2021-10-05 13:40:48.315686",
    "2021-10-05 13:40:48.315689",
    "2021-10-05 13:40:48.315691",
    "2021-10-05 13:40:48.315694",
    "2021-10-05 13:40:48.315696",
    "2021-10-05 13:40:48.315699",
    "2021-10-05 13:40:48.315701",
    "2021-10-05 13:40:48.315704",
    "2021-10-05 13:40:48.315706",
    "2021-10-05 13:40:48.315709",
    "2021-10-05 13:40:48.315711",
    "2021-10-05 13:40:48.315714",
    "2021-10-05 13:40:48.315716",
    "2021-10-05 13:40:48.315719",
    "2021-10-05 13:40:48.315722",
    "2021-10-05 13:40:48.315724",
    "2021-10-05 13:40:48.315726",
    "2021-10-05 13:40:48.315729",
    "2021-10-05 13:40:48.315731",
    "2021-10-05 13:40:48.315734",
    "2021-10-05 13:40:48.315736",
    "2021-10-05 13:40:48.315739",
    "2021-10-05 13:40:48.315742",
    "2021-10-05 13:40:48.315744",
    "2021-10-05 13:40:48.315747",
    "2021-10-05 13:40:48.315749",
    "2021-10-05 13:40:48.315752",
    "2021-10-05 13:40:48.315754",
    "2021-10-05 13:40:48.315756",
    "2021-10-05 13:40:48.315759",
    "2021-10-05 13:40:48.315764",
    "2021-10-05 13:40:48.315766",
    "2021-10-05 13:40:48.315769",
    "2021-10-05 13:40:48.315771",
    "2021-10-05 13:40:48.315774",
    "2021-10-05 13:40:48.315776",
    "2021-10-05 13:40:48.315779",
    "2021-10-05 13:40:48.315782",
    "2021-10-05 13:40:48.315784",
    "2021-10-05 13:40:48.315787",
    "2021-10-05 13:40:48.315789",
    "2021-10-05 13:40:48.315792",
    "2021-10-05 13:40:48.315794",
    "2021-10-05 13:40:48.315796",
    "2021-10-05 13:40:48.315799",
    "2021-10-05 13:40:48.315801",
    "2021-10-05 13:40:48.315804",
    "2021-10-05 13:40:48.315806",
    "2021-10-05 13:40:48.315809",
    "2021-10-05 13:40:48.315811",
    "2021-10-05 13:40:48.315814",
    "2021-10-05 13:40:48.315817",
    "2021-10-05 13:40:48.315819",
    "2021-10-05 13:40:48.315822",
    "2021-10-05 13:40:48.315824",
    "2021-10-05 13:40:48.315826",
    "2021-10-05 13:40:48.315829",
    "2021-10-05 13:40:48.315831",
    "2021-10-05 13:40:48.315834",
    "2021-10-05 13:40:48.315837"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

